$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force General/"@" text formatting for the assignment, then restore the
# default "Normal" style so the cell keeps its original (unstyled) look.
$numericLooking = @("D5", "D9", "D11", "D12", "D14", "D15", "D16", "D18", "D22", "D25", "D26", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D48", "D50")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.356.50'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '1.856.61'
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').Value = '314.26'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('E7').Value = '  -0.87%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '0.07324'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D11').Value = '19.93'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '0.07790'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = '1.887.33'
$ws.Range('E13').Value = '  +4.50%  '
$ws.Range('D14').Value = '5.384'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').Value = '6.546'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '91.87'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').Value = '0.000009057'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').Value = '27.369.77'
$ws.Range('E21').Value = '  +1.86%  '
$ws.Range('D22').Value = '5.135'
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = '2.127.48'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '1.928'
$ws.Range('E25').Value = '  +5.30%  '
$ws.Range('D26').Value = '152.14'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('E28').Value = '  -1.11%  '
$ws.Range('D29').Value = '5.107'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').Value = '116.02'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '0.08861'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').Value = '0.7699'
$ws.Range('E32').Value = '  +5.60%  '
$ws.Range('D33').Value = '3.039'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('E34').Value = '  +3.50%  '
$ws.Range('D35').Value = '4.501'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('D36').Value = '2.656'
$ws.Range('D37').Value = '0.01961'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '1.077'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '0.05233'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '2.957'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').Value = '7.026'
$ws.Range('E41').Value = '  -3.29%  '
$ws.Range('D42').Value = '0.5146'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').Value = '8.434'
$ws.Range('E44').Value = '  +2.53%  '
$ws.Range('D45').Value = '0.4823'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('D46').Value = '10.33'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('E47').Value = '  -0.70%  '
$ws.Range('D48').Value = '103.10'
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('D50').Value = '0.06221'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  +2.35%  '

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
